$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column G
$ws.Range("G1").Value = "RecoveredRatio"

# G2 was entered individually, so it is its own (non-shared) formula cell.
$ws.Range("G2").Formula = "=D2/B2"

# G3:G57 were entered/filled as one range, so Excel stores them as a single
# shared-formula group with G3 as the master cell (matches the diff).
$ws.Range("G3:G57").Formula = "=D3/B3"

# Select G2:G57 with G2 as the active cell, matching the post-edit selection.
$ws.Range("G2:G57").Select()
